$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1614.6923
$ws.Range("I2").Value = 2395.4443
$ws.Range("J2").Value = 1201.3529
$ws.Range("K2").Value = 2395.4443
$ws.Range("L2").Value = 1201.3529
$ws.Range("M2").Value = -2282.4443
$ws.Range("N2").Value = -1427.3529

$ws.Range("H4").Value = 1262.8
$ws.Range("I4").Value = 272
$ws.Range("J4").Value = 2749
$ws.Range("K4").Value = 272
$ws.Range("L4").Value = 2749
$ws.Range("M4").Value = -158
$ws.Range("N4").Value = -2977

$ws.Range("H17").Value = 977.9231
$ws.Range("J17").Value = 977.04
$ws.Range("L17").Value = 2931.12
$ws.Range("N17").Value = -3267.12

$ws.Range("H52").Value = 1000
$ws.Range("J52").Value = 1000
$ws.Range("L52").Value = 3000
$ws.Range("N52").Value = -3320

$ws.Range("H55").Value = 202.92857
$ws.Range("I55").Value = 195.53847
$ws.Range("J55").Value = 299
$ws.Range("K55").Value = 195.53847
$ws.Range("L55").Value = 299
$ws.Range("M55").Value = 18.46153000000001
$ws.Range("N55").Value = -727

$ws.Range("H58").Value = 5931.154
$ws.Range("J58").Value = 7410.5
$ws.Range("L58").Value = 22231.5
$ws.Range("N58").Value = -22531.5

$ws.Range("H106").Value = 37058470
$ws.Range("I106").Value = 43497164
$ws.Range("J106").Value = 36000
$ws.Range("K106").Value = 43497164
$ws.Range("L106").Value = 36000
$ws.Range("M106").Value = -43496533
$ws.Range("N106").Value = -37262

$ws.Range("H132").Value = 3626
$ws.Range("I132").Value = 1110
$ws.Range("K132").Value = 3330
$ws.Range("M132").Value = -800

$ws.Range("H138").Value = 4583.6
$ws.Range("I138").Value = 2704.6
$ws.Range("J138").Value = 5711
$ws.Range("K138").Value = 8113.799999999999
$ws.Range("L138").Value = 17133
$ws.Range("M138").Value = -2973.799999999999
$ws.Range("N138").Value = -27413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6453.769
$ws.Range("I63").Value = 1474.75
$ws.Range("J63").Value = 8666.666999999999
$ws.Range("K63").Value = 1474.75
$ws.Range("L63").Value = 8666.666999999999
$ws.Range("M63").Value = -788.75
$ws.Range("N63").Value = -10038.667

$ws.Range("H66").Value = 6453.769
$ws.Range("I66").Value = 1474.75
$ws.Range("J66").Value = 8666.666999999999
$ws.Range("K66").Value = 7373.75
$ws.Range("L66").Value = 43333.335
$ws.Range("M66").Value = -3941.75
$ws.Range("N66").Value = -50197.335

$ws.Range("H97").Value = 426.93332
$ws.Range("I97").Value = 481.27274
$ws.Range("J97").Value = 277.5
$ws.Range("K97").Value = 481.27274
$ws.Range("L97").Value = 277.5
$ws.Range("M97").Value = 14.72726
$ws.Range("N97").Value = -1269.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1412.1428
$ws.Range("I86").Value = 1305.7273
$ws.Range("J86").Value = 1802.3334
$ws.Range("K86").Value = 1305.7273
$ws.Range("L86").Value = 1802.3334
$ws.Range("M86").Value = -182.7273
$ws.Range("N86").Value = -4048.3334

$ws.Range("H89").Value = 1412.1428
$ws.Range("I89").Value = 1305.7273
$ws.Range("J89").Value = 1802.3334
$ws.Range("K89").Value = 6528.636500000001
$ws.Range("L89").Value = 9011.666999999999
$ws.Range("M89").Value = -912.6365000000005
$ws.Range("N89").Value = -20243.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2527
$ws.Range("I58").Value = 1897
$ws.Range("J58").Value = 4187.909
$ws.Range("K58").Value = 1897
$ws.Range("L58").Value = 4187.909
$ws.Range("M58").Value = -1694
$ws.Range("N58").Value = -4593.909

$ws.Range("H86").Value = 8466.333000000001
$ws.Range("I86").Value = 8466.333000000001
$ws.Range("K86").Value = 8466.333000000001
$ws.Range("M86").Value = -7343.333000000001

$ws.Range("H89").Value = 8466.333000000001
$ws.Range("I89").Value = 8466.333000000001
$ws.Range("K89").Value = 42331.665
$ws.Range("M89").Value = -36715.665

$ws.Range("H107").Value = 20834494
$ws.Range("I107").Value = 125000810
$ws.Range("K107").Value = 125000810
$ws.Range("M107").Value = -124998890

$ws.Range("H136").Value = 2527
$ws.Range("I136").Value = 1897
$ws.Range("J136").Value = 4187.909
$ws.Range("K136").Value = 5691
$ws.Range("L136").Value = 12563.727
$ws.Range("M136").Value = -3141
$ws.Range("N136").Value = -17663.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 511.72726
$ws.Range("I5").Value = 527.1429000000001
$ws.Range("J5").Value = 484.75
$ws.Range("K5").Value = 1581.4287
$ws.Range("L5").Value = 1454.25
$ws.Range("M5").Value = -1469.4287
$ws.Range("N5").Value = -1678.25

$ws.Range("H23").Value = 1564.25
$ws.Range("I23").Value = 562.5
$ws.Range("J23").Value = 2566
$ws.Range("K23").Value = 1687.5
$ws.Range("L23").Value = 7698
$ws.Range("M23").Value = -1452.5
$ws.Range("N23").Value = -8168

$ws.Range("H46").Value = 1668988.6
$ws.Range("I46").Value = 2200
$ws.Range("J46").Value = 3335777.2
$ws.Range("K46").Value = 6600
$ws.Range("L46").Value = 10007331.6
$ws.Range("M46").Value = -6509
$ws.Range("N46").Value = -10007513.6

$ws.Range("H56").Value = 12450.526
$ws.Range("I56").Value = 12450.526
$ws.Range("K56").Value = 12450.526
$ws.Range("M56").Value = -11920.526

$ws.Range("H68").Value = 1673.1538
$ws.Range("J68").Value = 1479.25
$ws.Range("L68").Value = 4437.75
$ws.Range("N68").Value = -6059.75

$ws.Range("H71").Value = 1673.1538
$ws.Range("J71").Value = 1479.25
$ws.Range("L71").Value = 13313.25
$ws.Range("N71").Value = -21425.25

$ws.Range("H88").Value = 2704.6
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 2704.6
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H132").Value = 4714.4287
$ws.Range("I132").Value = 2668
$ws.Range("J132").Value = 6249.25
$ws.Range("K132").Value = 24012
$ws.Range("L132").Value = 56243.25
$ws.Range("M132").Value = -21482
$ws.Range("N132").Value = -61303.25

$ws.Range("H133").Value = 3676.6667
$ws.Range("I133").Value = 3015
$ws.Range("K133").Value = 9045
$ws.Range("M133").Value = -3985

$ws.Range("H135").Value = 511.72726
$ws.Range("I135").Value = 527.1429000000001
$ws.Range("J135").Value = 484.75
$ws.Range("K135").Value = 4744.2861
$ws.Range("L135").Value = 4362.75
$ws.Range("M135").Value = -2209.2861
$ws.Range("N135").Value = -9432.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 299.0625
$ws.Range("I2").Value = 65.833336
$ws.Range("J2").Value = 439
$ws.Range("K2").Value = 65.833336
$ws.Range("L2").Value = 439
$ws.Range("M2").Value = 47.166664
$ws.Range("N2").Value = -665

$ws.Range("H7").Value = 11666.667
$ws.Range("J7").Value = 12500
$ws.Range("L7").Value = 12500
$ws.Range("N7").Value = -12724

$ws.Range("H8").Value = 11666.667
$ws.Range("J8").Value = 12500
$ws.Range("L8").Value = 12500
$ws.Range("N8").Value = -12778

$ws.Range("H11").Value = 6005081.5
$ws.Range("I11").Value = 8752622
$ws.Range("J11").Value = 510000
$ws.Range("K11").Value = 8752622
$ws.Range("L11").Value = 510000
$ws.Range("M11").Value = -8752483
$ws.Range("N11").Value = -510278

$ws.Range("H12").Value = 3007.1428
$ws.Range("I12").Value = 3060.6
$ws.Range("J12").Value = 2873.5
$ws.Range("K12").Value = 3060.6
$ws.Range("L12").Value = 2873.5
$ws.Range("M12").Value = -2920.6
$ws.Range("N12").Value = -3153.5

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H122").Value = 79323.92
$ws.Range("I122").Value = 2601
$ws.Range("K122").Value = 7803
$ws.Range("M122").Value = -5353

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 450
$ws.Range("J4").Value = 450
$ws.Range("L4").Value = 450
$ws.Range("N4").Value = -676

$ws.Range("H26").Value = 1500
$ws.Range("I26").Value = 1500
$ws.Range("K26").Value = 1500
$ws.Range("M26").Value = -1205

$ws.Range("H28").Value = 450
$ws.Range("J28").Value = 450
$ws.Range("L28").Value = 450
$ws.Range("N28").Value = -914

$ws.Range("H34").Value = 11000
$ws.Range("I34").Value = 13333.333
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 13333.333
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -13161.333
$ws.Range("N34").Value = -4344

$ws.Range("H37").Value = 450
$ws.Range("J37").Value = 450
$ws.Range("L37").Value = 450
$ws.Range("N37").Value = -664

$ws.Range("H46").Value = 1449.5714
$ws.Range("I46").Value = 749.4
$ws.Range("J46").Value = 3200
$ws.Range("K46").Value = 749.4
$ws.Range("L46").Value = 3200
$ws.Range("M46").Value = -561.4
$ws.Range("N46").Value = -3576

$ws.Range("H93").Value = 7842.2856
$ws.Range("I93").Value = 8255.111000000001
$ws.Range("K93").Value = 8255.111000000001
$ws.Range("M93").Value = -7007.111000000001

$ws.Range("H136").Value = 2887.4375
$ws.Range("I136").Value = 2585.8572
$ws.Range("K136").Value = 7757.571599999999
$ws.Range("M136").Value = -5207.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4605.875
$ws.Range("I132").Value = 2416.0557
$ws.Range("K132").Value = 7248.1671
$ws.Range("M132").Value = -4718.1671
